$wb = $excel.ActiveWorkbook

# This script applies numeric "want-to-go count" (F column) updates and
# a couple of "min ticket price" (G column) status changes to "已停售"
# (sold out / stopped selling), matching a data refresh across the four
# sheets: 展览, 演出, 本地生活, 全部类型.

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 2923
$ws.Range("F7").Value = 238
$ws.Range("F9").Value = 295
$ws.Range("F10").Value = 6830
$ws.Range("F11").Value = 35
$ws.Range("F12").Value = 45
$ws.Range("F13").Value = 342
$ws.Range("F14").Value = 596
$ws.Range("F15").Value = 1482
$ws.Range("F16").Value = 1107
$ws.Range("F17").Value = 2216
$ws.Range("F18").Value = 1461
$ws.Range("F19").Value = 647
$ws.Range("F21").Value = 1098
$ws.Range("F22").Value = 114
$ws.Range("F25").Value = 1679
$ws.Range("F26").Value = 1666
$ws.Range("F28").Value = 1027
$ws.Range("F31").Value = 1198
$ws.Range("F32").Value = 134
$ws.Range("F33").Value = 578
$ws.Range("F36").Value = 407
$ws.Range("F37").Value = 2
$ws.Range("F38").Value = 2441
$ws.Range("F39").Value = 2700
$ws.Range("F44").Value = 20
$ws.Range("F45").Value = 314
$ws.Range("F48").Value = 147

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 158
$ws.Range("F12").Value = 177
$ws.Range("F13").Value = 2
$ws.Range("F14").Value = 58
$ws.Range("F15").Value = 56
$ws.Range("F20").Value = 37
$ws.Range("F23").Value = 461
$ws.Range("F35").Value = 18
$ws.Range("F38").Value = 10

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 79
$ws.Range("F4").Value = 536
$ws.Range("F5").Value = 22
$ws.Range("G6").Value = "已停售"
$ws.Range("F7").Value = 1660
$ws.Range("F9").Value = 2705
$ws.Range("F10").Value = 994
$ws.Range("F11").Value = 898
$ws.Range("F13").Value = 245
$ws.Range("F14").Value = 1389
$ws.Range("F15").Value = 7270

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 536
$ws.Range("G4").Value = "已停售"
$ws.Range("F6").Value = 2923
$ws.Range("F7").Value = 238
$ws.Range("F8").Value = 1660
$ws.Range("F9").Value = 295
$ws.Range("F10").Value = 2705
$ws.Range("F11").Value = 6830
$ws.Range("F12").Value = 994
$ws.Range("F13").Value = 898
$ws.Range("F14").Value = 35
$ws.Range("F15").Value = 342
$ws.Range("F16").Value = 158
$ws.Range("F17").Value = 245
$ws.Range("F18").Value = 1389
$ws.Range("F19").Value = 596
$ws.Range("F20").Value = 2216
$ws.Range("F21").Value = 1461
$ws.Range("F22").Value = 647
$ws.Range("F24").Value = 1098
$ws.Range("F25").Value = 114
$ws.Range("F27").Value = 58
$ws.Range("F28").Value = 1679
$ws.Range("F29").Value = 1027
$ws.Range("F33").Value = 1198
$ws.Range("F36").Value = 461
$ws.Range("F37").Value = 407
$ws.Range("F39").Value = 2441
$ws.Range("F40").Value = 2700
$ws.Range("F44").Value = 314
